{"js": "// Update the East Asian / complex-script fonts recorded in the document's\n// paragraph styles (word/styles.xml), matching the upstream vignette's\n// switch from \"DejaVu Sans\" to \"Tahoma\" for the eastAsia font on the\n// styles that carried an explicit font, and recording the inherited\n// complex-script (\"DejaVu Sans\") font explicitly on the styles that\n// previously left it implicit.\n\nconst styles = context.document.getStyles();\n\n// Normal + Heading: eastAsia font DejaVu Sans -> Tahoma\nconst normal = styles.getByName(\"Normal\");\nnormal.font.nameFarEast = \"Tahoma\";\n\nconst heading = styles.getByName(\"Heading\");\nheading.font.nameFarEast = \"Tahoma\";\n\n// List, Caption, Index: record the complex-script (cs) font \"DejaVu Sans\"\n// explicitly on the style's own run properties.\nconst list = styles.getByName(\"List\");\nlist.font.nameBidirectional = \"DejaVu Sans\";\n\nconst caption = styles.getByName(\"Caption\");\ncaption.font.nameBidirectional = \"DejaVu Sans\";\n\nconst index = styles.getByName(\"Index\");\nindex.font.nameBidirectional = \"DejaVu Sans\";\n\nawait context.sync();\n", "ps1": "# Update the East Asian / complex-script fonts recorded in the document's\n# paragraph styles (word/styles.xml), matching the upstream vignette's\n# switch from \"DejaVu Sans\" to \"Tahoma\" for w:eastAsia on the styles that\n# carried an explicit font, and recording the complex-script font\n# (\"DejaVu Sans\") explicitly on the styles that previously inherited it.\n\n$d = $word.ActiveDocument\n\n# Normal + Heading: eastAsia font DejaVu Sans -> Tahoma\n$normal = $d.Styles(\"Normal\")\n$normal.Font.NameFarEast = \"Tahoma\"\n\n$heading = $d.Styles(\"Heading\")\n$heading.Font.NameFarEast = \"Tahoma\"\n\n# List, Caption, Index: record the inherited complex-script (cs) font\n# \"DejaVu Sans\" explicitly on the style's own run properties.\n$list = $d.Styles(\"List\")\n$list.Font.NameBi = \"DejaVu Sans\"\n\n$caption = $d.Styles(\"Caption\")\n$caption.Font.NameBi = \"DejaVu Sans\"\n\n$index = $d.Styles(\"Index\")\n$index.Font.NameBi = \"DejaVu Sans\"\n"}
